# BurningKnight3.xlsx edit:
#  - Rename sheets "Long Method" -> "Long_Method" and "Large Class" -> "Large_Class"
#  - Update the self-referencing label in cell B2 of each sheet to match the new name
#  - Update the B2 dropdown (data validation) list to use the underscore variants
#  - Switch the active sheet/tab from "Long Method" to "Large Class"
#  - Update the remembered selection on each sheet

$wb = $excel.ActiveWorkbook

$wsLongMethod = $wb.Worksheets.Item(1)
$wsLargeClass = $wb.Worksheets.Item(2)

# Rename worksheets (space -> underscore)
$wsLongMethod.Name = "Long_Method"
$wsLargeClass.Name = "Large_Class"

# Update the self-identifying text in B2 on each sheet
$wsLongMethod.Range("B2").Value = "Long_Method"
$wsLargeClass.Range("B2").Value = "Large_Class"

# Update the data validation dropdown list on B2 of each sheet
$wsLongMethod.Range("B2").Validation.Modify(3, 1, 1, '"Large_Class, Long_Method"')
$wsLargeClass.Range("B2").Validation.Modify(3, 1, 1, '"Large_Class, Long_Method"')

# Update remembered cell selection on each sheet
[void]$wsLongMethod.Range("B12").Select()

# Make "Large_Class" the active sheet/tab, with B2 selected
[void]$wsLargeClass.Activate()
[void]$wsLargeClass.Range("B2").Select()
